$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price string would otherwise be auto-parsed as a number
# by Excel; format as Text first so the literal string from the feed is kept.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the refreshed coin data (prices, 1h volume deltas, and the two
# ranking swaps at rows 47/48 and 50/51)
$ws.Range('D2').Value = '29.193.96'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.836.68'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '241.18'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = '0.6657'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.07375'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').Value = '0.2928'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D10').Value = '22.67'
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('D11').Value = '0.07726'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = '1.833.58'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '4.985'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = '0.6691'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').Value = '82.97'
$ws.Range('E15').Value = '  -5.03%  '
$ws.Range('D16').Value = '6.126'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '29.158.63'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '0.000008270'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('D19').Value = '225.50'
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('D20').Value = '12.46'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '7.130'
$ws.Range('E22').Value = '  -3.20%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '160.78'
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').Value = '8.625'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('D26').Value = '0.1392'
$ws.Range('E26').Value = '  -3.36%  '
$ws.Range('D27').Value = '17.97'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').Value = '1.510'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '4.113'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').Value = '4.034'
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').Value = '0.05305'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').Value = '0.7532'
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('D35').Value = '1.131'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').Value = '2.677'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('D37').Value = '1.298.10'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D39').Value = '2.721'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '0.9196'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').Value = '0.08738'
$ws.Range('E41').Value = '  +18.01%  '
$ws.Range('D42').Value = '5.961'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = '1.007'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').Value = '102.33'
$ws.Range('E44').Value = '  -2.29%  '
$ws.Range('D45').Value = '1.970.00'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '1.766'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000121'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('D49').Value = '63.16'
$ws.Range('E49').Value = '  -2.53%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05927'
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '8.984'
$ws.Range('E51').Value = '  -5.33%  '
